$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The two "Requisitos" entries (rows 24 and 25, columns B and C) swap their
# text content: row 24 now holds the LOT2052 entry, row 25 now holds the
# LOT2028 entry.

$lot2028 = "LOT2028 -  Tecnologia de Processos Fermentativos  (Requisito fraco)" + [char]10
$lot2052 = "LOT2052 -  Tecnologia de Bebidas Experimental  (Indicação de Conjunto)" + [char]10

$ws.Range("B24").Value = $lot2052
$ws.Range("C24").Value = $lot2052

$ws.Range("B25").Value = $lot2028
$ws.Range("C25").Value = $lot2028
